$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new ticker ("AVGO") at row 8, pushing the existing tickers
# (rows 8-74) down by one row (to rows 9-75), without disturbing the
# unrelated formula cell further down the sheet (row 77, "=+C78").
# Walk bottom-to-top so each value is copied before it gets overwritten.
for ($r = 74; $r -ge 8; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}

$ws.Cells.Item(8, 1).Value = "AVGO"

# Match the selection left behind in the saved file.
$ws.Range("A1:A75").Select() | Out-Null
